$d = $word.ActiveDocument

function Add-BoldSuffix {
    param(
        [string]$HeadingText,
        [string]$Suffix
    )
    $rng = $d.Content
    $found = $rng.Find.Execute($HeadingText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND: $HeadingText"
        return
    }
    $rng.Collapse(0)
    $rng.InsertAfter($Suffix)
    $rng.Font.Bold = $true
}

# 2. Contextualização - Massari
Add-BoldSuffix "2. Contextualização" " - Massari"

# 3. Definição de FDD - Massari
Add-BoldSuffix "3. Definição de FDD" " - Massari"

# 4. Princípios e Características do FDD - Ramon
Add-BoldSuffix "4. Princípios e Características do FDD" " - Ramon"

# 5. Processo do FDD - Ramon
Add-BoldSuffix "5. Processo do FDD" " - Ramon"

# 6. Papéis no FDD – Peuroca (en dash)
Add-BoldSuffix "6. Papéis no FDD" " – Peuroca"

# 7. Vantagens e Desvantagens - Peuroca
Add-BoldSuffix "7. Vantagens e Desvantagens" " - Peuroca"

# 8. Exemplo Prático / Estudo de Caso - Lucas
Add-BoldSuffix "8. Exemplo Prático / Estudo de Caso" " - Lucas"

# 9. Comparação com Outras Metodologias Ágeis - Lucas
Add-BoldSuffix "9. Comparação com Outras Metodologias Ágeis" " - Lucas"

# 10. Conclusão - Jean
Add-BoldSuffix "10. Conclusão" " - Jean"
